# Auto-generated script: apply row data swaps/rotations per commit diff
# Re-assigns id-row (column A) untouched; swaps/rotates columns B:AC between paired/grouped match rows
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- rows [19, 20] ---
# row 19 takes the original content of row 20
$ws.Range("B19").Value = 6815422
$ws.Range("C19").Value = 'Montenegro Prva Liga'
$ws.Range("D19").Value = 'Montenegro Prva Liga'
$ws.Range("E19").Value = 45151.625
$ws.Range("F19").Value = 'OFK Mladost DG'
$ws.Range("G19").Value = 'FK Decic Tuzi'
$ws.Range("H19").Value = 0
$ws.Range("I19").Value = 2
$ws.Range("J19").Value = 'A'
$ws.Range("K19").Value = 2.4
$ws.Range("L19").Value = 3
$ws.Range("M19").Value = 2.75
$ws.Range("N19").Value = 3.1
$ws.Range("O19").Value = 3
$ws.Range("P19").Value = 2.15
$ws.Range("Q19").Value = 0.25
$ws.Range("R19").Value = 1.875
$ws.Range("S19").Value = 1.925
$ws.Range("T19").Value = 2.25
$ws.Range("U19").Value = 2.025
$ws.Range("V19").Value = 1.775
$ws.Range("W19").Value = -1
$ws.Range("X19").Value = -1
$ws.Range("Y19").Value = 1.15
$ws.Range("Z19").Value = -1
$ws.Range("AA19").Value = 0.925
$ws.Range("AB19").Value = -0.5
$ws.Range("AC19").Value = 0.3875
# row 20 takes the original content of row 19
$ws.Range("B20").Value = 6815305
$ws.Range("C20").Value = 'Montenegro Prva Liga'
$ws.Range("D20").Value = 'Montenegro Prva Liga'
$ws.Range("E20").Value = 45151.625
$ws.Range("F20").Value = 'Buducnost Podgorica'
$ws.Range("G20").Value = 'FK Jezero'
$ws.Range("H20").Value = 1
$ws.Range("I20").Value = 1
$ws.Range("J20").Value = 'D'
$ws.Range("K20").Value = 1.3
$ws.Range("L20").Value = 5
$ws.Range("M20").Value = 7
$ws.Range("N20").Value = 1.571
$ws.Range("O20").Value = 4
$ws.Range("P20").Value = 4.2
$ws.Range("Q20").Value = -0.75
$ws.Range("R20").Value = 1.75
$ws.Range("S20").Value = 1.95
$ws.Range("T20").Value = 2.5
$ws.Range("U20").Value = 1.95
$ws.Range("V20").Value = 1.85
$ws.Range("W20").Value = -1
$ws.Range("X20").Value = 3
$ws.Range("Y20").Value = -1
$ws.Range("Z20").Value = -1
$ws.Range("AA20").Value = 0.95
$ws.Range("AB20").Value = -1
$ws.Range("AC20").Value = 0.8500000000000001

# --- rows [25, 26] ---
# row 25 takes the original content of row 26
$ws.Range("B25").Value = 6815308
$ws.Range("C25").Value = 'Montenegro Prva Liga'
$ws.Range("D25").Value = 'Montenegro Prva Liga'
$ws.Range("E25").Value = 45158.625
$ws.Range("F25").Value = 'FK Jezero'
$ws.Range("G25").Value = 'OFK Petrovac'
$ws.Range("H25").Value = 2
$ws.Range("I25").Value = 0
$ws.Range("J25").Value = 'H'
$ws.Range("K25").Value = 2.3
$ws.Range("L25").Value = 3
$ws.Range("M25").Value = 2.9
$ws.Range("N25").Value = 2.7
$ws.Range("O25").Value = 2.625
$ws.Range("P25").Value = 2.75
$ws.Range("Q25").Value = 0
$ws.Range("R25").Value = 1.875
$ws.Range("S25").Value = 1.925
$ws.Range("T25").Value = 2.25
$ws.Range("U25").Value = 1.925
$ws.Range("V25").Value = 1.875
$ws.Range("W25").Value = 1.7
$ws.Range("X25").Value = -1
$ws.Range("Y25").Value = -1
$ws.Range("Z25").Value = 0.875
$ws.Range("AA25").Value = -1
$ws.Range("AB25").Value = -0.5
$ws.Range("AC25").Value = 0.4375
# row 26 takes the original content of row 25
$ws.Range("B26").Value = 6815423
$ws.Range("C26").Value = 'Montenegro Prva Liga'
$ws.Range("D26").Value = 'Montenegro Prva Liga'
$ws.Range("E26").Value = 45158.625
$ws.Range("F26").Value = 'FK Decic Tuzi'
$ws.Range("G26").Value = 'FK Mornar Bar'
$ws.Range("H26").Value = 3
$ws.Range("I26").Value = 1
$ws.Range("J26").Value = 'H'
$ws.Range("K26").Value = 1.85
$ws.Range("L26").Value = 3.25
$ws.Range("M26").Value = 3.75
$ws.Range("N26").Value = 1.8
$ws.Range("O26").Value = 3.2
$ws.Range("P26").Value = 4.2
$ws.Range("Q26").Value = -0.5
$ws.Range("R26").Value = 1.825
$ws.Range("S26").Value = 1.975
$ws.Range("T26").Value = 2
$ws.Range("U26").Value = 1.825
$ws.Range("V26").Value = 1.975
$ws.Range("W26").Value = 0.8
$ws.Range("X26").Value = -1
$ws.Range("Y26").Value = -1
$ws.Range("Z26").Value = 0.825
$ws.Range("AA26").Value = -1
$ws.Range("AB26").Value = 0.825
$ws.Range("AC26").Value = -1

# --- rows [30, 31] ---
# row 30 takes the original content of row 31
$ws.Range("B30").Value = 6815315
$ws.Range("C30").Value = 'Montenegro Prva Liga'
$ws.Range("D30").Value = 'Montenegro Prva Liga'
$ws.Range("E30").Value = 45164.625
$ws.Range("F30").Value = 'FK Decic Tuzi'
$ws.Range("G30").Value = 'FK Rudar Pljevlja'
$ws.Range("H30").Value = 0
$ws.Range("I30").Value = 1
$ws.Range("J30").Value = 'A'
$ws.Range("K30").Value = 1.615
$ws.Range("L30").Value = 3.5
$ws.Range("M30").Value = 4.75
$ws.Range("N30").Value = 1.4
$ws.Range("O30").Value = 3.8
$ws.Range("P30").Value = 6.5
$ws.Range("Q30").Value = -1.25
$ws.Range("R30").Value = 2
$ws.Range("S30").Value = 1.8
$ws.Range("T30").Value = 2.5
$ws.Range("U30").Value = 1.95
$ws.Range("V30").Value = 1.85
$ws.Range("W30").Value = -1
$ws.Range("X30").Value = -1
$ws.Range("Y30").Value = 5.5
$ws.Range("Z30").Value = -1
$ws.Range("AA30").Value = 0.8
$ws.Range("AB30").Value = -1
$ws.Range("AC30").Value = 0.8500000000000001
# row 31 takes the original content of row 30
$ws.Range("B31").Value = 6815312
$ws.Range("C31").Value = 'Montenegro Prva Liga'
$ws.Range("D31").Value = 'Montenegro Prva Liga'
$ws.Range("E31").Value = 45164.625
$ws.Range("F31").Value = 'Buducnost Podgorica'
$ws.Range("G31").Value = 'FK Arsenal'
$ws.Range("H31").Value = 2
$ws.Range("I31").Value = 0
$ws.Range("J31").Value = 'H'
$ws.Range("K31").Value = 1.444
$ws.Range("L31").Value = 4
$ws.Range("M31").Value = 6
$ws.Range("N31").Value = 1.4
$ws.Range("O31").Value = 4
$ws.Range("P31").Value = 6.5
$ws.Range("Q31").Value = -1.25
$ws.Range("R31").Value = 1.95
$ws.Range("S31").Value = 1.85
$ws.Range("T31").Value = 2.5
$ws.Range("U31").Value = 1.775
$ws.Range("V31").Value = 1.925
$ws.Range("W31").Value = 0.3999999999999999
$ws.Range("X31").Value = -1
$ws.Range("Y31").Value = -1
$ws.Range("Z31").Value = 0.95
$ws.Range("AA31").Value = -1
$ws.Range("AB31").Value = -1
$ws.Range("AC31").Value = 0.925

# --- rows [38, 39] ---
# row 38 takes the original content of row 39
$ws.Range("B38").Value = 6815322
$ws.Range("C38").Value = 'Montenegro Prva Liga'
$ws.Range("D38").Value = 'Montenegro Prva Liga'
$ws.Range("E38").Value = 45185.60416666666
$ws.Range("F38").Value = 'OFK Mladost DG'
$ws.Range("G38").Value = 'FK Arsenal'
$ws.Range("H38").Value = 2
$ws.Range("I38").Value = 1
$ws.Range("J38").Value = 'H'
$ws.Range("K38").Value = 2.375
$ws.Range("L38").Value = 3
$ws.Range("M38").Value = 2.75
$ws.Range("N38").Value = 2.625
$ws.Range("O38").Value = 3
$ws.Range("P38").Value = 2.4
$ws.Range("Q38").Value = 0
$ws.Range("R38").Value = 2
$ws.Range("S38").Value = 1.8
$ws.Range("T38").Value = 2
$ws.Range("U38").Value = 1.725
$ws.Range("V38").Value = 1.975
$ws.Range("W38").Value = 1.625
$ws.Range("X38").Value = -1
$ws.Range("Y38").Value = -1
$ws.Range("Z38").Value = 1
$ws.Range("AA38").Value = -1
$ws.Range("AB38").Value = 0.7250000000000001
$ws.Range("AC38").Value = -1
# row 39 takes the original content of row 38
$ws.Range("B39").Value = 6815321
$ws.Range("C39").Value = 'Montenegro Prva Liga'
$ws.Range("D39").Value = 'Montenegro Prva Liga'
$ws.Range("E39").Value = 45185.60416666666
$ws.Range("F39").Value = 'OFK Petrovac'
$ws.Range("G39").Value = 'FK Jedinstvo Bijelo Polje'
$ws.Range("H39").Value = 1
$ws.Range("I39").Value = 1
$ws.Range("J39").Value = 'D'
$ws.Range("K39").Value = 1.8
$ws.Range("L39").Value = 3.4
$ws.Range("M39").Value = 3.75
$ws.Range("N39").Value = 1.6
$ws.Range("O39").Value = 3.5
$ws.Range("P39").Value = 4.75
$ws.Range("Q39").Value = -0.75
$ws.Range("R39").Value = 1.825
$ws.Range("S39").Value = 1.975
$ws.Range("T39").Value = 2.5
$ws.Range("U39").Value = 1.95
$ws.Range("V39").Value = 1.75
$ws.Range("W39").Value = -1
$ws.Range("X39").Value = 2.5
$ws.Range("Y39").Value = -1
$ws.Range("Z39").Value = -1
$ws.Range("AA39").Value = 0.9750000000000001
$ws.Range("AB39").Value = -1
$ws.Range("AC39").Value = 0.75

# --- rows [48, 49] ---
# row 48 takes the original content of row 49
$ws.Range("B48").Value = 6815331
$ws.Range("C48").Value = 'Montenegro Prva Liga'
$ws.Range("D48").Value = 'Montenegro Prva Liga'
$ws.Range("E48").Value = 45196.54166666666
$ws.Range("F48").Value = 'FK Jedinstvo Bijelo Polje'
$ws.Range("G48").Value = 'FK Decic Tuzi'
$ws.Range("H48").Value = 0
$ws.Range("I48").Value = 1
$ws.Range("J48").Value = 'A'
$ws.Range("K48").Value = 3
$ws.Range("L48").Value = 3.1
$ws.Range("M48").Value = 2.2
$ws.Range("N48").Value = 5.75
$ws.Range("O48").Value = 3.6
$ws.Range("P48").Value = 1.5
$ws.Range("Q48").Value = 1
$ws.Range("R48").Value = 1.825
$ws.Range("S48").Value = 1.975
$ws.Range("T48").Value = 2.25
$ws.Range("U48").Value = 1.95
$ws.Range("V48").Value = 1.85
$ws.Range("W48").Value = -1
$ws.Range("X48").Value = -1
$ws.Range("Y48").Value = 0.5
$ws.Range("Z48").Value = 0
$ws.Range("AA48").Value = -0
$ws.Range("AB48").Value = -1
$ws.Range("AC48").Value = 0.8500000000000001
# row 49 takes the original content of row 48
$ws.Range("B49").Value = 6815333
$ws.Range("C49").Value = 'Montenegro Prva Liga'
$ws.Range("D49").Value = 'Montenegro Prva Liga'
$ws.Range("E49").Value = 45196.54166666666
$ws.Range("F49").Value = 'Sutjeska Niksic'
$ws.Range("G49").Value = 'FK Jezero'
$ws.Range("H49").Value = 1
$ws.Range("I49").Value = 1
$ws.Range("J49").Value = 'D'
$ws.Range("K49").Value = 1.5
$ws.Range("L49").Value = 3.6
$ws.Range("M49").Value = 6
$ws.Range("N49").Value = 1.5
$ws.Range("O49").Value = 3.6
$ws.Range("P49").Value = 5.75
$ws.Range("Q49").Value = -1
$ws.Range("R49").Value = 1.9
$ws.Range("S49").Value = 1.9
$ws.Range("T49").Value = 2.25
$ws.Range("U49").Value = 1.85
$ws.Range("V49").Value = 1.95
$ws.Range("W49").Value = -1
$ws.Range("X49").Value = 2.6
$ws.Range("Y49").Value = -1
$ws.Range("Z49").Value = -1
$ws.Range("AA49").Value = 0.8999999999999999
$ws.Range("AB49").Value = -0.5
$ws.Range("AC49").Value = 0.475

# --- rows [53, 55] ---
# row 53 takes the original content of row 55
$ws.Range("B53").Value = 6815334
$ws.Range("C53").Value = 'Montenegro Prva Liga'
$ws.Range("D53").Value = 'Montenegro Prva Liga'
$ws.Range("E53").Value = 45200.54166666666
$ws.Range("F53").Value = 'Sutjeska Niksic'
$ws.Range("G53").Value = 'FK Mornar Bar'
$ws.Range("H53").Value = 0
$ws.Range("I53").Value = 1
$ws.Range("J53").Value = 'A'
$ws.Range("K53").Value = 1.444
$ws.Range("L53").Value = 4
$ws.Range("M53").Value = 6.5
$ws.Range("N53").Value = 1.444
$ws.Range("O53").Value = 4
$ws.Range("P53").Value = 6.5
$ws.Range("Q53").Value = -1.25
$ws.Range("R53").Value = 2
$ws.Range("S53").Value = 1.8
$ws.Range("T53").Value = 2.5
$ws.Range("U53").Value = 2
$ws.Range("V53").Value = 1.8
$ws.Range("W53").Value = -1
$ws.Range("X53").Value = -1
$ws.Range("Y53").Value = 5.5
$ws.Range("Z53").Value = -1
$ws.Range("AA53").Value = 0.8
$ws.Range("AB53").Value = -1
$ws.Range("AC53").Value = 0.8
# row 55 takes the original content of row 53
$ws.Range("B55").Value = 7279987
$ws.Range("C55").Value = 'Montenegro Prva Liga'
$ws.Range("D55").Value = 'Montenegro Prva Liga'
$ws.Range("E55").Value = 45200.54166666666
$ws.Range("F55").Value = 'FK Jezero'
$ws.Range("G55").Value = 'FK Arsenal'
$ws.Range("H55").Value = 1
$ws.Range("I55").Value = 1
$ws.Range("J55").Value = 'D'
$ws.Range("K55").Value = 2.1
$ws.Range("L55").Value = 3
$ws.Range("M55").Value = 3.25
$ws.Range("N55").Value = 2.05
$ws.Range("O55").Value = 3
$ws.Range("P55").Value = 3.4
$ws.Range("Q55").Value = -0.25
$ws.Range("R55").Value = 1.8
$ws.Range("S55").Value = 2
$ws.Range("T55").Value = 2
$ws.Range("U55").Value = 1.925
$ws.Range("V55").Value = 1.875
$ws.Range("W55").Value = -1
$ws.Range("X55").Value = 2
$ws.Range("Y55").Value = -1
$ws.Range("Z55").Value = -0.5
$ws.Range("AA55").Value = 0.5
$ws.Range("AB55").Value = 0
$ws.Range("AC55").Value = -0

# --- rows [62, 63, 64] ---
# row 62 takes the original content of row 63
$ws.Range("B62").Value = 7366684
$ws.Range("C62").Value = 'Montenegro Prva Liga'
$ws.Range("D62").Value = 'Montenegro Prva Liga'
$ws.Range("E62").Value = 45220.41666666666
$ws.Range("F62").Value = 'FK Rudar Pljevlja'
$ws.Range("G62").Value = 'OFK Petrovac'
$ws.Range("H62").Value = 1
$ws.Range("I62").Value = 0
$ws.Range("J62").Value = 'H'
$ws.Range("K62").Value = 2.875
$ws.Range("L62").Value = 2.9
$ws.Range("M62").Value = 2.375
$ws.Range("N62").Value = 2.625
$ws.Range("O62").Value = 2.9
$ws.Range("P62").Value = 2.55
$ws.Range("Q62").Value = 0
$ws.Range("R62").Value = 1.925
$ws.Range("S62").Value = 1.875
$ws.Range("T62").Value = 2.25
$ws.Range("U62").Value = 1.925
$ws.Range("V62").Value = 1.875
$ws.Range("W62").Value = 1.625
$ws.Range("X62").Value = -1
$ws.Range("Y62").Value = -1
$ws.Range("Z62").Value = 0.925
$ws.Range("AA62").Value = -1
$ws.Range("AB62").Value = -1
$ws.Range("AC62").Value = 0.875
# row 63 takes the original content of row 64
$ws.Range("B63").Value = 7366683
$ws.Range("C63").Value = 'Montenegro Prva Liga'
$ws.Range("D63").Value = 'Montenegro Prva Liga'
$ws.Range("E63").Value = 45220.41666666666
$ws.Range("F63").Value = 'FK Arsenal'
$ws.Range("G63").Value = 'FK Mornar Bar'
$ws.Range("H63").Value = 2
$ws.Range("I63").Value = 2
$ws.Range("J63").Value = 'D'
$ws.Range("K63").Value = 2.375
$ws.Range("L63").Value = 2.8
$ws.Range("M63").Value = 3
$ws.Range("N63").Value = 2.3
$ws.Range("O63").Value = 2.7
$ws.Range("P63").Value = 3.3
$ws.Range("Q63").Value = -0.25
$ws.Range("R63").Value = 2
$ws.Range("S63").Value = 1.8
$ws.Range("T63").Value = 1.75
$ws.Range("U63").Value = 1.875
$ws.Range("V63").Value = 1.925
$ws.Range("W63").Value = -1
$ws.Range("X63").Value = 1.7
$ws.Range("Y63").Value = -1
$ws.Range("Z63").Value = -0.5
$ws.Range("AA63").Value = 0.4
$ws.Range("AB63").Value = 0.875
$ws.Range("AC63").Value = -1
# row 64 takes the original content of row 62
$ws.Range("B64").Value = 6815343
$ws.Range("C64").Value = 'Montenegro Prva Liga'
$ws.Range("D64").Value = 'Montenegro Prva Liga'
$ws.Range("E64").Value = 45220.41666666666
$ws.Range("F64").Value = 'Sutjeska Niksic'
$ws.Range("G64").Value = 'FK Jedinstvo Bijelo Polje'
$ws.Range("H64").Value = 2
$ws.Range("I64").Value = 0
$ws.Range("J64").Value = 'H'
$ws.Range("K64").Value = 1.333
$ws.Range("L64").Value = 4.2
$ws.Range("M64").Value = 8
$ws.Range("N64").Value = 1.333
$ws.Range("O64").Value = 4.2
$ws.Range("P64").Value = 8
$ws.Range("Q64").Value = -1.5
$ws.Range("R64").Value = 1.975
$ws.Range("S64").Value = 1.825
$ws.Range("T64").Value = 2.75
$ws.Range("U64").Value = 1.9
$ws.Range("V64").Value = 1.9
$ws.Range("W64").Value = 0.333
$ws.Range("X64").Value = -1
$ws.Range("Y64").Value = -1
$ws.Range("Z64").Value = 0.9750000000000001
$ws.Range("AA64").Value = -1
$ws.Range("AB64").Value = -1
$ws.Range("AC64").Value = 0.8999999999999999

# --- rows [76, 77] ---
# row 76 takes the original content of row 77
$ws.Range("B76").Value = 6815357
$ws.Range("C76").Value = 'Montenegro Prva Liga'
$ws.Range("D76").Value = 'Montenegro Prva Liga'
$ws.Range("E76").Value = 45241.5
$ws.Range("F76").Value = 'OFK Mladost DG'
$ws.Range("G76").Value = 'Sutjeska Niksic'
$ws.Range("H76").Value = 1
$ws.Range("I76").Value = 1
$ws.Range("J76").Value = 'D'
$ws.Range("K76").Value = 4.8
$ws.Range("L76").Value = 3.5
$ws.Range("M76").Value = 1.615
$ws.Range("N76").Value = 4.75
$ws.Range("O76").Value = 3.6
$ws.Range("P76").Value = 1.571
$ws.Range("Q76").Value = 1
$ws.Range("R76").Value = 1.775
$ws.Range("S76").Value = 2.025
$ws.Range("T76").Value = 2.25
$ws.Range("U76").Value = 1.825
$ws.Range("V76").Value = 1.975
$ws.Range("W76").Value = -1
$ws.Range("X76").Value = 2.6
$ws.Range("Y76").Value = -1
$ws.Range("Z76").Value = 0.7749999999999999
$ws.Range("AA76").Value = -1
$ws.Range("AB76").Value = -0.5
$ws.Range("AC76").Value = 0.4875
# row 77 takes the original content of row 76
$ws.Range("B77").Value = 6815358
$ws.Range("C77").Value = 'Montenegro Prva Liga'
$ws.Range("D77").Value = 'Montenegro Prva Liga'
$ws.Range("E77").Value = 45241.5
$ws.Range("F77").Value = 'OFK Petrovac'
$ws.Range("G77").Value = 'FK Arsenal'
$ws.Range("H77").Value = 1
$ws.Range("I77").Value = 1
$ws.Range("J77").Value = 'D'
$ws.Range("K77").Value = 2.1
$ws.Range("L77").Value = 3.1
$ws.Range("M77").Value = 3.2
$ws.Range("N77").Value = 1.75
$ws.Range("O77").Value = 3.3
$ws.Range("P77").Value = 4.2
$ws.Range("Q77").Value = -0.5
$ws.Range("R77").Value = 1.8
$ws.Range("S77").Value = 2
$ws.Range("T77").Value = 2.25
$ws.Range("U77").Value = 1.95
$ws.Range("V77").Value = 1.85
$ws.Range("W77").Value = -1
$ws.Range("X77").Value = 2.3
$ws.Range("Y77").Value = -1
$ws.Range("Z77").Value = -1
$ws.Range("AA77").Value = 1
$ws.Range("AB77").Value = -0.5
$ws.Range("AC77").Value = 0.425

# --- rows [81, 82] ---
# row 81 takes the original content of row 82
$ws.Range("B81").Value = 6815362
$ws.Range("C81").Value = 'Montenegro Prva Liga'
$ws.Range("D81").Value = 'Montenegro Prva Liga'
$ws.Range("E81").Value = 45256.45833333334
$ws.Range("F81").Value = 'Sutjeska Niksic'
$ws.Range("G81").Value = 'FK Decic Tuzi'
$ws.Range("H81").Value = 1
$ws.Range("I81").Value = 1
$ws.Range("J81").Value = 'D'
$ws.Range("K81").Value = 2.2
$ws.Range("L81").Value = 3
$ws.Range("M81").Value = 3.1
$ws.Range("N81").Value = 2.375
$ws.Range("O81").Value = 2.875
$ws.Range("P81").Value = 3
$ws.Range("Q81").Value = -0.25
$ws.Range("R81").Value = 2.05
$ws.Range("S81").Value = 1.75
$ws.Range("T81").Value = 2
$ws.Range("U81").Value = 1.8
$ws.Range("V81").Value = 2
$ws.Range("W81").Value = -1
$ws.Range("X81").Value = 1.875
$ws.Range("Y81").Value = -1
$ws.Range("Z81").Value = -0.5
$ws.Range("AA81").Value = 0.375
$ws.Range("AB81").Value = 0
$ws.Range("AC81").Value = -0
# row 82 takes the original content of row 81
$ws.Range("B82").Value = 6815430
$ws.Range("C82").Value = 'Montenegro Prva Liga'
$ws.Range("D82").Value = 'Montenegro Prva Liga'
$ws.Range("E82").Value = 45256.45833333334
$ws.Range("F82").Value = 'Buducnost Podgorica'
$ws.Range("G82").Value = 'FK Mornar Bar'
$ws.Range("H82").Value = 4
$ws.Range("I82").Value = 3
$ws.Range("J82").Value = 'H'
$ws.Range("K82").Value = 1.444
$ws.Range("L82").Value = 3.75
$ws.Range("M82").Value = 6.5
$ws.Range("N82").Value = 1.4
$ws.Range("O82").Value = 4
$ws.Range("P82").Value = 7
$ws.Range("Q82").Value = -1.25
$ws.Range("R82").Value = 1.875
$ws.Range("S82").Value = 1.925
$ws.Range("T82").Value = 2.5
$ws.Range("U82").Value = 1.775
$ws.Range("V82").Value = 1.925
$ws.Range("W82").Value = 0.3999999999999999
$ws.Range("X82").Value = -1
$ws.Range("Y82").Value = -1
$ws.Range("Z82").Value = -0.5
$ws.Range("AA82").Value = 0.4625
$ws.Range("AB82").Value = 0.7749999999999999
$ws.Range("AC82").Value = -1

# --- rows [85, 86] ---
# row 85 takes the original content of row 86
$ws.Range("B85").Value = 6815365
$ws.Range("C85").Value = 'Montenegro Prva Liga'
$ws.Range("D85").Value = 'Montenegro Prva Liga'
$ws.Range("E85").Value = 45262.375
$ws.Range("F85").Value = 'FK Rudar Pljevlja'
$ws.Range("G85").Value = 'Sutjeska Niksic'
$ws.Range("H85").Value = 1
$ws.Range("I85").Value = 2
$ws.Range("J85").Value = 'A'
$ws.Range("K85").Value = 4.75
$ws.Range("L85").Value = 3.5
$ws.Range("M85").Value = 1.615
$ws.Range("N85").Value = 5
$ws.Range("O85").Value = 3.6
$ws.Range("P85").Value = 1.615
$ws.Range("Q85").Value = 1
$ws.Range("R85").Value = 1.725
$ws.Range("S85").Value = 1.975
$ws.Range("T85").Value = 2.25
$ws.Range("U85").Value = 1.875
$ws.Range("V85").Value = 1.925
$ws.Range("W85").Value = -1
$ws.Range("X85").Value = -1
$ws.Range("Y85").Value = 0.615
$ws.Range("Z85").Value = 0
$ws.Range("AA85").Value = -0
$ws.Range("AB85").Value = 0.875
$ws.Range("AC85").Value = -1
# row 86 takes the original content of row 85
$ws.Range("B86").Value = 6815366
$ws.Range("C86").Value = 'Montenegro Prva Liga'
$ws.Range("D86").Value = 'Montenegro Prva Liga'
$ws.Range("E86").Value = 45262.375
$ws.Range("F86").Value = 'FK Decic Tuzi'
$ws.Range("G86").Value = 'FK Arsenal'
$ws.Range("H86").Value = 0
$ws.Range("I86").Value = 1
$ws.Range("J86").Value = 'A'
$ws.Range("K86").Value = 1.571
$ws.Range("L86").Value = 3.6
$ws.Range("M86").Value = 5
$ws.Range("N86").Value = 1.363
$ws.Range("O86").Value = 4.2
$ws.Range("P86").Value = 7
$ws.Range("Q86").Value = -1.25
$ws.Range("R86").Value = 1.825
$ws.Range("S86").Value = 1.975
$ws.Range("T86").Value = 2.75
$ws.Range("U86").Value = 1.975
$ws.Range("V86").Value = 1.825
$ws.Range("W86").Value = -1
$ws.Range("X86").Value = -1
$ws.Range("Y86").Value = 6
$ws.Range("Z86").Value = -1
$ws.Range("AA86").Value = 0.9750000000000001
$ws.Range("AB86").Value = -1
$ws.Range("AC86").Value = 0.825

# --- rows [100, 101] ---
# row 100 takes the original content of row 101
$ws.Range("B100").Value = 6815433
$ws.Range("C100").Value = 'Montenegro Prva Liga'
$ws.Range("D100").Value = 'Montenegro Prva Liga'
$ws.Range("E100").Value = 45346.41666666666
$ws.Range("F100").Value = 'OFK Mladost DG'
$ws.Range("G100").Value = 'FK Mornar Bar'
$ws.Range("H100").Value = 1
$ws.Range("I100").Value = 2
$ws.Range("J100").Value = 'A'
$ws.Range("K100").Value = 2.5
$ws.Range("L100").Value = 3
$ws.Range("M100").Value = 2.6
$ws.Range("N100").Value = 2.5
$ws.Range("O100").Value = 3
$ws.Range("P100").Value = 2.6
$ws.Range("Q100").Value = 0
$ws.Range("R100").Value = 1.85
$ws.Range("S100").Value = 1.95
$ws.Range("T100").Value = 2
$ws.Range("U100").Value = 1.975
$ws.Range("V100").Value = 1.825
$ws.Range("W100").Value = -1
$ws.Range("X100").Value = -1
$ws.Range("Y100").Value = 1.6
$ws.Range("Z100").Value = -1
$ws.Range("AA100").Value = 0.95
$ws.Range("AB100").Value = 0.9750000000000001
$ws.Range("AC100").Value = -1
# row 101 takes the original content of row 100
$ws.Range("B101").Value = 6815378
$ws.Range("C101").Value = 'Montenegro Prva Liga'
$ws.Range("D101").Value = 'Montenegro Prva Liga'
$ws.Range("E101").Value = 45346.41666666666
$ws.Range("F101").Value = 'FK Rudar Pljevlja'
$ws.Range("G101").Value = 'Buducnost Podgorica'
$ws.Range("H101").Value = 3
$ws.Range("I101").Value = 0
$ws.Range("J101").Value = 'H'
$ws.Range("K101").Value = 10
$ws.Range("L101").Value = 5.5
$ws.Range("M101").Value = 1.2
$ws.Range("N101").Value = 11
$ws.Range("O101").Value = 4.75
$ws.Range("P101").Value = 1.222
$ws.Range("Q101").Value = 1.75
$ws.Range("R101").Value = 1.85
$ws.Range("S101").Value = 1.95
$ws.Range("T101").Value = 2.5
$ws.Range("U101").Value = 1.825
$ws.Range("V101").Value = 1.975
$ws.Range("W101").Value = 10
$ws.Range("X101").Value = -1
$ws.Range("Y101").Value = -1
$ws.Range("Z101").Value = 0.8500000000000001
$ws.Range("AA101").Value = -1
$ws.Range("AB101").Value = 0.825
$ws.Range("AC101").Value = -1

# --- rows [107, 108] ---
# row 107 takes the original content of row 108
$ws.Range("B107").Value = 7890508
$ws.Range("C107").Value = 'Montenegro Prva Liga'
$ws.Range("D107").Value = 'Montenegro Prva Liga'
$ws.Range("E107").Value = 45350.5625
$ws.Range("F107").Value = 'OFK Petrovac'
$ws.Range("G107").Value = 'FK Rudar Pljevlja'
$ws.Range("H107").Value = 1
$ws.Range("I107").Value = 1
$ws.Range("J107").Value = 'D'
$ws.Range("K107").Value = 1.75
$ws.Range("L107").Value = 3.1
$ws.Range("M107").Value = 4.5
$ws.Range("N107").Value = 1.8
$ws.Range("O107").Value = 3.2
$ws.Range("P107").Value = 4
$ws.Range("Q107").Value = -0.5
$ws.Range("R107").Value = 1.875
$ws.Range("S107").Value = 1.925
$ws.Range("T107").Value = 2.25
$ws.Range("U107").Value = 1.95
$ws.Range("V107").Value = 1.85
$ws.Range("W107").Value = -1
$ws.Range("X107").Value = 2.2
$ws.Range("Y107").Value = -1
$ws.Range("Z107").Value = -1
$ws.Range("AA107").Value = 0.925
$ws.Range("AB107").Value = -0.5
$ws.Range("AC107").Value = 0.425
# row 108 takes the original content of row 107
$ws.Range("B108").Value = 7890506
$ws.Range("C108").Value = 'Montenegro Prva Liga'
$ws.Range("D108").Value = 'Montenegro Prva Liga'
$ws.Range("E108").Value = 45350.5625
$ws.Range("F108").Value = 'FK Mornar Bar'
$ws.Range("G108").Value = 'FK Arsenal'
$ws.Range("H108").Value = 0
$ws.Range("I108").Value = 0
$ws.Range("J108").Value = 'D'
$ws.Range("K108").Value = 1.85
$ws.Range("L108").Value = 3.1
$ws.Range("M108").Value = 3.9
$ws.Range("N108").Value = 1.85
$ws.Range("O108").Value = 3.3
$ws.Range("P108").Value = 3.5
$ws.Range("Q108").Value = -0.5
$ws.Range("R108").Value = 1.925
$ws.Range("S108").Value = 1.875
$ws.Range("T108").Value = 2
$ws.Range("U108").Value = 1.95
$ws.Range("V108").Value = 1.85
$ws.Range("W108").Value = -1
$ws.Range("X108").Value = 2.3
$ws.Range("Y108").Value = -1
$ws.Range("Z108").Value = -1
$ws.Range("AA108").Value = 0.875
$ws.Range("AB108").Value = -1
$ws.Range("AC108").Value = 0.8500000000000001

# --- rows [117, 118] ---
# row 117 takes the original content of row 118
$ws.Range("B117").Value = 6815389
$ws.Range("C117").Value = 'Montenegro Prva Liga'
$ws.Range("D117").Value = 'Montenegro Prva Liga'
$ws.Range("E117").Value = 45360.58333333334
$ws.Range("F117").Value = 'FK Mornar Bar'
$ws.Range("G117").Value = 'FK Jedinstvo Bijelo Polje'
$ws.Range("H117").Value = 0
$ws.Range("I117").Value = 0
$ws.Range("J117").Value = 'D'
$ws.Range("K117").Value = 1.909
$ws.Range("L117").Value = 3.1
$ws.Range("M117").Value = 3.75
$ws.Range("N117").Value = 1.75
$ws.Range("O117").Value = 3.2
$ws.Range("P117").Value = 4.333
$ws.Range("Q117").Value = -0.5
$ws.Range("R117").Value = 1.8
$ws.Range("S117").Value = 2
$ws.Range("T117").Value = 2
$ws.Range("U117").Value = 1.9
$ws.Range("V117").Value = 1.9
$ws.Range("W117").Value = -1
$ws.Range("X117").Value = 2.2
$ws.Range("Y117").Value = -1
$ws.Range("Z117").Value = -1
$ws.Range("AA117").Value = 1
$ws.Range("AB117").Value = -1
$ws.Range("AC117").Value = 0.8999999999999999
# row 118 takes the original content of row 117
$ws.Range("B118").Value = 6815393
$ws.Range("C118").Value = 'Montenegro Prva Liga'
$ws.Range("D118").Value = 'Montenegro Prva Liga'
$ws.Range("E118").Value = 45360.58333333334
$ws.Range("F118").Value = 'FK Decic Tuzi'
$ws.Range("G118").Value = 'FK Rudar Pljevlja'
$ws.Range("H118").Value = 3
$ws.Range("I118").Value = 1
$ws.Range("J118").Value = 'H'
$ws.Range("K118").Value = 1.285
$ws.Range("L118").Value = 4.75
$ws.Range("M118").Value = 8
$ws.Range("N118").Value = 1.333
$ws.Range("O118").Value = 4.5
$ws.Range("P118").Value = 7
$ws.Range("Q118").Value = -1.25
$ws.Range("R118").Value = 1.75
$ws.Range("S118").Value = 1.95
$ws.Range("T118").Value = 2.25
$ws.Range("U118").Value = 1.75
$ws.Range("V118").Value = 1.95
$ws.Range("W118").Value = 0.333
$ws.Range("X118").Value = -1
$ws.Range("Y118").Value = -1
$ws.Range("Z118").Value = 0.75
$ws.Range("AA118").Value = -1
$ws.Range("AB118").Value = 0.75
$ws.Range("AC118").Value = -1

# --- rows [126, 127] ---
# row 126 takes the original content of row 127
$ws.Range("B126").Value = 6815401
$ws.Range("C126").Value = 'Montenegro Prva Liga'
$ws.Range("D126").Value = 'Montenegro Prva Liga'
$ws.Range("E126").Value = 45368.45833333334
$ws.Range("F126").Value = 'FK Decic Tuzi'
$ws.Range("G126").Value = 'Sutjeska Niksic'
$ws.Range("H126").Value = 0
$ws.Range("I126").Value = 0
$ws.Range("J126").Value = 'D'
$ws.Range("K126").Value = 2.55
$ws.Range("L126").Value = 3
$ws.Range("M126").Value = 2.6
$ws.Range("N126").Value = 2.1
$ws.Range("O126").Value = 3.1
$ws.Range("P126").Value = 3.3
$ws.Range("Q126").Value = -0.25
$ws.Range("R126").Value = 1.825
$ws.Range("S126").Value = 1.975
$ws.Range("T126").Value = 2
$ws.Range("U126").Value = 1.925
$ws.Range("V126").Value = 1.875
$ws.Range("W126").Value = -1
$ws.Range("X126").Value = 2.1
$ws.Range("Y126").Value = -1
$ws.Range("Z126").Value = -0.5
$ws.Range("AA126").Value = 0.4875
$ws.Range("AB126").Value = -1
$ws.Range("AC126").Value = 0.875
# row 127 takes the original content of row 126
$ws.Range("B127").Value = 6815402
$ws.Range("C127").Value = 'Montenegro Prva Liga'
$ws.Range("D127").Value = 'Montenegro Prva Liga'
$ws.Range("E127").Value = 45368.45833333334
$ws.Range("F127").Value = 'FK Rudar Pljevlja'
$ws.Range("G127").Value = 'FK Jezero'
$ws.Range("H127").Value = 0
$ws.Range("I127").Value = 1
$ws.Range("J127").Value = 'A'
$ws.Range("K127").Value = 2.8
$ws.Range("L127").Value = 3
$ws.Range("M127").Value = 2.375
$ws.Range("N127").Value = 2.45
$ws.Range("O127").Value = 2.9
$ws.Range("P127").Value = 2.75
$ws.Range("Q127").Value = 0
$ws.Range("R127").Value = 1.775
$ws.Range("S127").Value = 2.025
$ws.Range("T127").Value = 1.75
$ws.Range("U127").Value = 1.825
$ws.Range("V127").Value = 1.975
$ws.Range("W127").Value = -1
$ws.Range("X127").Value = -1
$ws.Range("Y127").Value = 1.75
$ws.Range("Z127").Value = -1
$ws.Range("AA127").Value = 1.025
$ws.Range("AB127").Value = -1
$ws.Range("AC127").Value = 0.9750000000000001

# --- rows [130, 131] ---
# row 130 takes the original content of row 131
$ws.Range("B130").Value = 6815403
$ws.Range("C130").Value = 'Montenegro Prva Liga'
$ws.Range("D130").Value = 'Montenegro Prva Liga'
$ws.Range("E130").Value = 45381.45833333334
$ws.Range("F130").Value = 'FK Jezero'
$ws.Range("G130").Value = 'FK Mornar Bar'
$ws.Range("H130").Value = 3
$ws.Range("I130").Value = 1
$ws.Range("J130").Value = 'H'
$ws.Range("K130").Value = 2.875
$ws.Range("L130").Value = 2.75
$ws.Range("M130").Value = 2.5
$ws.Range("N130").Value = 3.1
$ws.Range("O130").Value = 2.55
$ws.Range("P130").Value = 2.55
$ws.Range("Q130").Value = 0
$ws.Range("R130").Value = 2.05
$ws.Range("S130").Value = 1.75
$ws.Range("T130").Value = 1.75
$ws.Range("U130").Value = 1.95
$ws.Range("V130").Value = 1.85
$ws.Range("W130").Value = 2.1
$ws.Range("X130").Value = -1
$ws.Range("Y130").Value = -1
$ws.Range("Z130").Value = 1.05
$ws.Range("AA130").Value = -1
$ws.Range("AB130").Value = 0.95
$ws.Range("AC130").Value = -1
# row 131 takes the original content of row 130
$ws.Range("B131").Value = 6815404
$ws.Range("C131").Value = 'Montenegro Prva Liga'
$ws.Range("D131").Value = 'Montenegro Prva Liga'
$ws.Range("E131").Value = 45381.45833333334
$ws.Range("F131").Value = 'Sutjeska Niksic'
$ws.Range("G131").Value = 'FK Rudar Pljevlja'
$ws.Range("H131").Value = 2
$ws.Range("I131").Value = 1
$ws.Range("J131").Value = 'H'
$ws.Range("K131").Value = 1.5
$ws.Range("L131").Value = 3.75
$ws.Range("M131").Value = 5.75
$ws.Range("N131").Value = 1.5
$ws.Range("O131").Value = 3.75
$ws.Range("P131").Value = 5.75
$ws.Range("Q131").Value = -1
$ws.Range("R131").Value = 1.85
$ws.Range("S131").Value = 1.95
$ws.Range("T131").Value = 2.25
$ws.Range("U131").Value = 1.9
$ws.Range("V131").Value = 1.9
$ws.Range("W131").Value = 0.5
$ws.Range("X131").Value = -1
$ws.Range("Y131").Value = -1
$ws.Range("Z131").Value = 0
$ws.Range("AA131").Value = -0
$ws.Range("AB131").Value = 0.8999999999999999
$ws.Range("AC131").Value = -1

# --- rows [135, 136] ---
# row 135 takes the original content of row 136
$ws.Range("B135").Value = 8043517
$ws.Range("C135").Value = 'Montenegro Prva Liga'
$ws.Range("D135").Value = 'Montenegro Prva Liga'
$ws.Range("E135").Value = 45385.41666666666
$ws.Range("F135").Value = 'FK Jedinstvo Bijelo Polje'
$ws.Range("G135").Value = 'FK Decic Tuzi'
$ws.Range("H135").Value = 0
$ws.Range("I135").Value = 2
$ws.Range("J135").Value = 'A'
$ws.Range("K135").Value = 5.5
$ws.Range("L135").Value = 3.2
$ws.Range("M135").Value = 1.615
$ws.Range("N135").Value = 6.5
$ws.Range("O135").Value = 3.4
$ws.Range("P135").Value = 1.533
$ws.Range("Q135").Value = 1
$ws.Range("R135").Value = 1.85
$ws.Range("S135").Value = 1.95
$ws.Range("T135").Value = 2.25
$ws.Range("U135").Value = 2.025
$ws.Range("V135").Value = 1.775
$ws.Range("W135").Value = -1
$ws.Range("X135").Value = -1
$ws.Range("Y135").Value = 0.5329999999999999
$ws.Range("Z135").Value = -1
$ws.Range("AA135").Value = 0.95
$ws.Range("AB135").Value = -0.5
$ws.Range("AC135").Value = 0.3875
# row 136 takes the original content of row 135
$ws.Range("B136").Value = 8043518
$ws.Range("C136").Value = 'Montenegro Prva Liga'
$ws.Range("D136").Value = 'Montenegro Prva Liga'
$ws.Range("E136").Value = 45385.41666666666
$ws.Range("F136").Value = 'FK Arsenal'
$ws.Range("G136").Value = 'FK Rudar Pljevlja'
$ws.Range("H136").Value = 4
$ws.Range("I136").Value = 2
$ws.Range("J136").Value = 'H'
$ws.Range("K136").Value = 1.909
$ws.Range("L136").Value = 3
$ws.Range("M136").Value = 3.9
$ws.Range("N136").Value = 1.65
$ws.Range("O136").Value = 3.3
$ws.Range("P136").Value = 5
$ws.Range("Q136").Value = -0.75
$ws.Range("R136").Value = 1.875
$ws.Range("S136").Value = 1.925
$ws.Range("T136").Value = 2
$ws.Range("U136").Value = 1.8
$ws.Range("V136").Value = 2
$ws.Range("W136").Value = 0.6499999999999999
$ws.Range("X136").Value = -1
$ws.Range("Y136").Value = -1
$ws.Range("Z136").Value = 0.875
$ws.Range("AA136").Value = -1
$ws.Range("AB136").Value = 0.8
$ws.Range("AC136").Value = -1

# --- rows [141, 142] ---
# row 141 takes the original content of row 142
$ws.Range("B141").Value = 8062094
$ws.Range("C141").Value = 'Montenegro Prva Liga'
$ws.Range("D141").Value = 'Montenegro Prva Liga'
$ws.Range("E141").Value = 45389.41666666666
$ws.Range("F141").Value = 'FK Rudar Pljevlja'
$ws.Range("G141").Value = 'FK Jedinstvo Bijelo Polje'
$ws.Range("H141").Value = 1
$ws.Range("I141").Value = 1
$ws.Range("J141").Value = 'D'
$ws.Range("K141").Value = 2.25
$ws.Range("L141").Value = 3
$ws.Range("M141").Value = 3
$ws.Range("N141").Value = 2.25
$ws.Range("O141").Value = 3.1
$ws.Range("P141").Value = 2.875
$ws.Range("Q141").Value = -0.25
$ws.Range("R141").Value = 2
$ws.Range("S141").Value = 1.8
$ws.Range("T141").Value = 2.25
$ws.Range("U141").Value = 1.95
$ws.Range("V141").Value = 1.85
$ws.Range("W141").Value = -1
$ws.Range("X141").Value = 2.1
$ws.Range("Y141").Value = -1
$ws.Range("Z141").Value = -0.5
$ws.Range("AA141").Value = 0.4
$ws.Range("AB141").Value = -0.5
$ws.Range("AC141").Value = 0.425
# row 142 takes the original content of row 141
$ws.Range("B142").Value = 8062093
$ws.Range("C142").Value = 'Montenegro Prva Liga'
$ws.Range("D142").Value = 'Montenegro Prva Liga'
$ws.Range("E142").Value = 45389.41666666666
$ws.Range("F142").Value = 'FK Jezero'
$ws.Range("G142").Value = 'FK Arsenal'
$ws.Range("H142").Value = 4
$ws.Range("I142").Value = 0
$ws.Range("J142").Value = 'H'
$ws.Range("K142").Value = 2.1
$ws.Range("L142").Value = 3
$ws.Range("M142").Value = 3.25
$ws.Range("N142").Value = 2.1
$ws.Range("O142").Value = 3
$ws.Range("P142").Value = 3.2
$ws.Range("Q142").Value = -0.25
$ws.Range("R142").Value = 1.875
$ws.Range("S142").Value = 1.925
$ws.Range("T142").Value = 2.25
$ws.Range("U142").Value = 1.95
$ws.Range("V142").Value = 1.85
$ws.Range("W142").Value = 1.1
$ws.Range("X142").Value = -1
$ws.Range("Y142").Value = -1
$ws.Range("Z142").Value = 0.875
$ws.Range("AA142").Value = -1
$ws.Range("AB142").Value = 0.95
$ws.Range("AC142").Value = -1

